{"js": "// Apply the \"story map\" edit: trim the FR/NF/B/K legend paragraph down to\n// just FR + NF, drop the now-empty paragraph that followed it, rewrite the\n// functional/non-functional requirement rows into full user-story text, and\n// remove the rows that the author dropped from the tables (FR-15, NF-05,\n// NF-06, NF-07).\n\n// ---------------------------------------------------------------------\n// 1. Trim the legend paragraph: remove the \"B stands for limitation\" /\n//    \"K stands for quality requirement.\" runs (and their leading line\n//    break), then delete the empty paragraph that used to follow it.\n// ---------------------------------------------------------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet legendParagraph = null;\nlet emptyParagraphAfterLegend = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"stands for non-functional requirement.\") !== -1 &&\n      t.indexOf(\"stands for limitation\") !== -1) {\n    legendParagraph = paragraphs.items[i];\n    emptyParagraphAfterLegend = paragraphs.items[i + 1];\n    break;\n  }\n}\n\nif (legendParagraph) {\n  const tailSearch = legendParagraph.search(\n    \"non-functional requirement.\\u000bB stands for limitation\\u000bK stands for quality requirement.\",\n    { matchCase: true }\n  );\n  tailSearch.load(\"items\");\n  await context.sync();\n  if (tailSearch.items.length > 0) {\n    tailSearch.items[0].insertText(\"non-functional requirement.\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\nif (emptyParagraphAfterLegend) {\n  emptyParagraphAfterLegend.delete();\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2. Helper: set the plain-text content of a table cell's first\n//    paragraph while preserving its paragraph-level formatting.\n// ---------------------------------------------------------------------\nasync function setCellText(cell, newText) {\n  cell.body.paragraphs.load(\"items\");\n  await context.sync();\n  const para = cell.body.paragraphs.items[0];\n  para.getRange().insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3. Functional requirements table (first table).\n// ---------------------------------------------------------------------\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst functionalTable = tables.items[0];\nconst functionalRows = functionalTable.rows;\nfunctionalRows.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < functionalRows.items.length; i++) {\n  functionalRows.items[i].cells.load(\"items\");\n}\nawait context.sync();\n\nconst frTextByRowIndex = {\n  1: \"As a user I want to be able to select categories so I can find a potential desired product(s).\",\n  2: \"As a user I want to be able to add one or multiple items to my shopping cart so I can order my desired products.\",\n  4: \"As a user, I want to be able to create an account so that my data is saved for a subsequent order.\",\n  5: \"As an administrator, I want to be able to easily create user and administrator accounts so that selecting roles remains under administrators.\",\n  6: \"As an administrator, I want to be able to modify products in the web shop so that the information and price remain up-to-date.\",\n  7: \"As a user, I want to be able to check out the products in my cart so that the order is initiated.\",\n  8: \"As a user, I want to be able to update my account information so that my information can be kept up-to-date.\",\n  9: \"As a user, I want to be able to search for products all over the site so that I can find my desired product.\",\n  10: \"As a user, I want to be able to filter by certain characteristics in the selected category so that I can more easily find the product I am looking for.\",\n  11: \"As a user, I want to see related products on the product page so that I can order something with them if necessary.\",\n  12: \"As a developer, I want a product to be able to contain multiple categories so that searching for particular products is simplified.\",\n  13: \"As a user, I want to be able to log in so I can update my account information.\",\n  14: \"As a user, I want to be able to register so I can update my account information.\"\n};\n\nfor (const rowIndex of Object.keys(frTextByRowIndex)) {\n  const idx = Number(rowIndex);\n  const textCell = functionalRows.items[idx].cells.items[1];\n  await setCellText(textCell, frTextByRowIndex[idx]);\n}\n\n// Row 15 (FR-15 \"The project will be unit tested\") was removed entirely.\nfunctionalRows.items[15].delete();\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 4. Non-functional requirements table (second table).\n// ---------------------------------------------------------------------\nconst tables2 = context.document.body.tables;\ntables2.load(\"items\");\nawait context.sync();\n\nconst nonFunctionalTable = tables2.items[1];\nconst nfRows = nonFunctionalTable.rows;\nnfRows.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < nfRows.items.length; i++) {\n  nfRows.items[i].cells.load(\"items\");\n}\nawait context.sync();\n\nconst nfTextByRowIndex = {\n  1: \"As a user I want to be able to view the website on every device whenever I want so I can purchase products everywhere.\",\n  2: \"As a software developer I want the application to be maintainable and scalable so I can add functionality to it in the future.\",\n  3: \"As a software developer I want the front-end and backend communicating via an API.\",\n  4: \"As a user I want the website to load within 1 seconds so I won\\u2019t visit another website.\"\n};\n\nfor (const rowIndex of Object.keys(nfTextByRowIndex)) {\n  const idx = Number(rowIndex);\n  const textCell = nfRows.items[idx].cells.items[1];\n  await setCellText(textCell, nfTextByRowIndex[idx]);\n}\n\n// Rows 5, 6, 7 (NF-05, NF-06, NF-07) were removed entirely. Delete from the\n// end so earlier indices stay valid.\nnfRows.items[7].delete();\nnfRows.items[6].delete();\nnfRows.items[5].delete();\nawait context.sync();\n", "ps1": "# Apply the \"story map\" edit: trim the FR/NF/B/K legend paragraph down to\n# just FR + NF, drop the now-empty paragraph that followed it, rewrite the\n# functional/non-functional requirement rows into full user-story text, and\n# remove the rows that the author dropped from the tables (FR-15, NF-05,\n# NF-06, NF-07).\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1. Trim the legend paragraph: remove the \"B stands for limitation\" /\n#    \"K stands for quality requirement.\" runs (and their leading line\n#    break) by replacing the whole tail with just the NF sentence.\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"non-functional requirement.\" + [char]11 + \"B stands for limitation\" + [char]11 + \"K stands for quality requirement.\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"non-functional requirement.\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# The paragraph that used to just hold an empty line after the legend is\n# now the 5th paragraph in the document; delete it outright.\n$legendGapParagraph = $d.Paragraphs.Item(5)\n$legendGapParagraph.Range.Delete()\n\n# ---------------------------------------------------------------------\n# 2. Functional requirements table (first table).\n# ---------------------------------------------------------------------\n$functionalTable = $d.Tables.Item(1)\n\n$functionalTable.Cell(2, 2).Range.Text = \"As a user I want to be able to select categories so I can find a potential desired product(s).\"\n$functionalTable.Cell(3, 2).Range.Text = \"As a user I want to be able to add one or multiple items to my shopping cart so I can order my desired products.\"\n$functionalTable.Cell(5, 2).Range.Text = \"As a user, I want to be able to create an account so that my data is saved for a subsequent order.\"\n$functionalTable.Cell(6, 2).Range.Text = \"As an administrator, I want to be able to easily create user and administrator accounts so that selecting roles remains under administrators.\"\n$functionalTable.Cell(7, 2).Range.Text = \"As an administrator, I want to be able to modify products in the web shop so that the information and price remain up-to-date.\"\n$functionalTable.Cell(8, 2).Range.Text = \"As a user, I want to be able to check out the products in my cart so that the order is initiated.\"\n$functionalTable.Cell(9, 2).Range.Text = \"As a user, I want to be able to update my account information so that my information can be kept up-to-date.\"\n$functionalTable.Cell(10, 2).Range.Text = \"As a user, I want to be able to search for products all over the site so that I can find my desired product.\"\n$functionalTable.Cell(11, 2).Range.Text = \"As a user, I want to be able to filter by certain characteristics in the selected category so that I can more easily find the product I am looking for.\"\n$functionalTable.Cell(12, 2).Range.Text = \"As a user, I want to see related products on the product page so that I can order something with them if necessary.\"\n$functionalTable.Cell(13, 2).Range.Text = \"As a developer, I want a product to be able to contain multiple categories so that searching for particular products is simplified.\"\n$functionalTable.Cell(14, 2).Range.Text = \"As a user, I want to be able to log in so I can update my account information.\"\n$functionalTable.Cell(15, 2).Range.Text = \"As a user, I want to be able to register so I can update my account information.\"\n\n# Row 16 (FR-15 \"The project will be unit tested\") was removed entirely.\n$functionalTable.Rows.Item(16).Delete()\n\n# ---------------------------------------------------------------------\n# 3. Non-functional requirements table (second table).\n# ---------------------------------------------------------------------\n$nonFunctionalTable = $d.Tables.Item(2)\n\n$nonFunctionalTable.Cell(2, 2).Range.Text = \"As a user I want to be able to view the website on every device whenever I want so I can purchase products everywhere.\"\n$nonFunctionalTable.Cell(3, 2).Range.Text = \"As a software developer I want the application to be maintainable and scalable so I can add functionality to it in the future.\"\n$nonFunctionalTable.Cell(4, 2).Range.Text = \"As a software developer I want the front-end and backend communicating via an API.\"\n$nonFunctionalTable.Cell(5, 2).Range.Text = \"As a user I want the website to load within 1 seconds so I won\" + [char]8217 + \"t visit another website.\"\n\n# Rows 6, 7, 8 (NF-05, NF-06, NF-07) were removed entirely. Delete from the\n# end so earlier row indices stay valid.\n$nonFunctionalTable.Rows.Item(8).Delete()\n$nonFunctionalTable.Rows.Item(7).Delete()\n$nonFunctionalTable.Rows.Item(6).Delete()\n"}
